$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (source old row 29)
$ws.Range("A4").Value = 74387725
$ws.Range("Q4").Value = 376996.7937296238
$ws.Range("R4").Value = 6700002.117591534

# Row 5 (source old row 4)
$ws.Range("A5").Value = 74387724
$ws.Range("Q5").Value = 377000.7918530595
$ws.Range("R5").Value = 6700265.926565454

# Row 6 (source old row 5)
$ws.Range("A6").Value = 74387739
$ws.Range("B6").Value = 77506
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = 'Garnlav'
$ws.Range("G6").Value = 'Alectoria sarmentosa'
$ws.Range("H6").Value = '(Ach.) Ach.'
$ws.Range("Q6").Value = 376872.8419193771
$ws.Range("R6").Value = 6700376.039685452

# Row 7 (source old row 30)
$ws.Range("A7").Value = 74387758
$ws.Range("Q7").Value = 376849.9079317059
$ws.Range("R7").Value = 6700079.764070438

# Row 8 (source old row 6)
$ws.Range("A8").Value = 74387736
$ws.Range("B8").Value = 77258
$ws.Range("E8").Value = 6446
$ws.Range("F8").Value = 'Kolflarnlav'
$ws.Range("G8").Value = 'Carbonicola anthracophila'
$ws.Range("H8").Value = '(Nyl.) Bendiksby & Timdal'
$ws.Range("Q8").Value = 376889.9729257529
$ws.Range("R8").Value = 6700268.204837496

# Row 9 (source old row 31)
$ws.Range("A9").Value = 74387721
$ws.Range("Q9").Value = 377025.0852994107
$ws.Range("R9").Value = 6699945.796444287

# Row 10 (source old row 7)
$ws.Range("A10").Value = 74387738
$ws.Range("Q10").Value = 376873.8805943133
$ws.Range("R10").Value = 6700260.843082374

# Row 11 (source old row 32)
$ws.Range("A11").Value = 74387719
$ws.Range("Q11").Value = 377043.0654386982
$ws.Range("R11").Value = 6699935.794955891

# Row 12 (source old row 8)
$ws.Range("A12").Value = 74387729
$ws.Range("Q12").Value = 376970.887149193
$ws.Range("R12").Value = 6700260.021504997

# Row 13 (source old row 9)
$ws.Range("A13").Value = 74387747
$ws.Range("B13").Value = 77506
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = 'Garnlav'
$ws.Range("G13").Value = 'Alectoria sarmentosa'
$ws.Range("H13").Value = '(Ach.) Ach.'
$ws.Range("Q13").Value = 376867.0591878957
$ws.Range("R13").Value = 6700409.844704853

# Row 14 (source old row 10)
$ws.Range("A14").Value = 74387731
$ws.Range("Q14").Value = 376967.8946852534
$ws.Range("R14").Value = 6700273.962341819

# Row 15 (source old row 33)
$ws.Range("A15").Value = 74387754
$ws.Range("Q15").Value = 376856.2017830843
$ws.Range("R15").Value = 6699929.778228731

# Row 16 (source old row 11)
$ws.Range("A16").Value = 74387730
$ws.Range("Q16").Value = 376969.1101584906
$ws.Range("R16").Value = 6700266.012934987

# Row 17 (source old row 12)
$ws.Range("A17").Value = 74387743
$ws.Range("Q17").Value = 376870.2103046329
$ws.Range("R17").Value = 6700386.013980007

# Row 18 (source old row 34)
$ws.Range("A18").Value = 74387720
$ws.Range("Q18").Value = 377037.1978320775
$ws.Range("R18").Value = 6699937.97113596

# Row 19 (source old row 13)
$ws.Range("A19").Value = 74387759
$ws.Range("B19").Value = 81236
$ws.Range("E19").Value = 1312
$ws.Range("F19").Value = 'Gammelgransskål'
$ws.Range("G19").Value = 'Pseudographis pinicola'
$ws.Range("H19").Value = '(Nyl.) Rehm'
$ws.Range("Q19").Value = 376838.8620604862
$ws.Range("R19").Value = 6700191.846452983

# Row 20 (source old row 14)
$ws.Range("A20").Value = 74387735
$ws.Range("Q20").Value = 376899.2007959604
$ws.Range("R20").Value = 6700292.110244696

# Row 21 (source old row 35)
$ws.Range("A21").Value = 74387749
$ws.Range("B21").Value = 77506
$ws.Range("E21").Value = 6425
$ws.Range("F21").Value = 'Garnlav'
$ws.Range("G21").Value = 'Alectoria sarmentosa'
$ws.Range("H21").Value = '(Ach.) Ach.'
$ws.Range("Q21").Value = 376864.8522645697
$ws.Range("R21").Value = 6699951.233937764

# Row 22 (source old row 15)
$ws.Range("A22").Value = 74387726
$ws.Range("B22").Value = 77506
$ws.Range("P22").Value = 'Jonsmyrsudden, Vrm'
$ws.Range("Q22").Value = 376979.8799874667
$ws.Range("R22").Value = 6700277.015484395
$ws.Range("Y22").NumberFormat = "@"
$ws.Range("Y22").Value = '2018-11-15'
$ws.Range("Z22").NumberFormat = "@"
$ws.Range("Z22").Value = '00:00'
$ws.Range("AA22").NumberFormat = "@"
$ws.Range("AA22").Value = '2018-11-15'
$ws.Range("AB22").NumberFormat = "@"
$ws.Range("AB22").Value = '00:00'
$ws.Range("AW22").Value = 'Per Gustafsson'
$ws.Range("AX22").Value = 'Per Gustafsson'

# Row 23 (source old row 16)
$ws.Range("A23").Value = 74387752
$ws.Range("B23").Value = 77506
$ws.Range("E23").Value = 6425
$ws.Range("F23").Value = 'Garnlav'
$ws.Range("G23").Value = 'Alectoria sarmentosa'
$ws.Range("H23").Value = '(Ach.) Ach.'
$ws.Range("P23").Value = 'Jonsmyrsudden, Vrm'
$ws.Range("Q23").Value = 376863.8981990077
$ws.Range("R23").Value = 6700243.882853701
$ws.Range("Y23").NumberFormat = "@"
$ws.Range("Y23").Value = '2018-11-15'
$ws.Range("Z23").NumberFormat = "@"
$ws.Range("Z23").Value = '00:00'
$ws.Range("AA23").NumberFormat = "@"
$ws.Range("AA23").Value = '2018-11-15'
$ws.Range("AB23").NumberFormat = "@"
$ws.Range("AB23").Value = '00:00'
$ws.Range("AW23").Value = 'Per Gustafsson'
$ws.Range("AX23").Value = 'Per Gustafsson'

# Row 24 (source old row 17)
$ws.Range("A24").Value = 74387757
$ws.Range("B24").Value = 77506
$ws.Range("P24").Value = 'Jonsmyrsudden, Vrm'
$ws.Range("Q24").Value = 376851.2419536079
$ws.Range("R24").Value = 6700221.082168487
$ws.Range("Y24").NumberFormat = "@"
$ws.Range("Y24").Value = '2018-11-15'
$ws.Range("Z24").NumberFormat = "@"
$ws.Range("Z24").Value = '00:00'
$ws.Range("AA24").NumberFormat = "@"
$ws.Range("AA24").Value = '2018-11-15'
$ws.Range("AB24").NumberFormat = "@"
$ws.Range("AB24").Value = '00:00'
$ws.Range("AW24").Value = 'Per Gustafsson'
$ws.Range("AX24").Value = 'Per Gustafsson'

# Row 25 (source old row 18)
$ws.Range("A25").Value = 74387728
$ws.Range("B25").Value = 77506
$ws.Range("E25").Value = 6425
$ws.Range("F25").Value = 'Garnlav'
$ws.Range("G25").Value = 'Alectoria sarmentosa'
$ws.Range("H25").Value = '(Ach.) Ach.'
$ws.Range("P25").Value = 'Jonsmyrsudden, Vrm'
$ws.Range("Q25").Value = 376971.0292153726
$ws.Range("R25").Value = 6700278.798573329
$ws.Range("Y25").NumberFormat = "@"
$ws.Range("Y25").Value = '2018-11-15'
$ws.Range("Z25").NumberFormat = "@"
$ws.Range("Z25").Value = '00:00'
$ws.Range("AA25").NumberFormat = "@"
$ws.Range("AA25").Value = '2018-11-15'
$ws.Range("AB25").NumberFormat = "@"
$ws.Range("AB25").Value = '00:00'
$ws.Range("AW25").Value = 'Per Gustafsson'
$ws.Range("AX25").Value = 'Per Gustafsson'

# Row 26 (source old row 36)
$ws.Range("A26").Value = 74387742
$ws.Range("B26").Value = 77506
$ws.Range("P26").Value = 'Jonsmyrsudden, Vrm'
$ws.Range("Q26").Value = 376869.9907624553
$ws.Range("R26").Value = 6700058.81623072
$ws.Range("Y26").NumberFormat = "@"
$ws.Range("Y26").Value = '2018-11-15'
$ws.Range("Z26").NumberFormat = "@"
$ws.Range("Z26").Value = '00:00'
$ws.Range("AA26").NumberFormat = "@"
$ws.Range("AA26").Value = '2018-11-15'
$ws.Range("AB26").NumberFormat = "@"
$ws.Range("AB26").Value = '00:00'
$ws.Range("AW26").Value = 'Per Gustafsson'
$ws.Range("AX26").Value = 'Per Gustafsson'

# Row 27 (source old row 19)
$ws.Range("A27").Value = 74387760
$ws.Range("B27").Value = 77506
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = 'Garnlav'
$ws.Range("G27").Value = 'Alectoria sarmentosa'
$ws.Range("H27").Value = '(Ach.) Ach.'
$ws.Range("P27").Value = 'Jonsmyrsudden, Vrm'
$ws.Range("Q27").Value = 376837.7977669895
$ws.Range("R27").Value = 6700204.239387355
$ws.Range("Y27").NumberFormat = "@"
$ws.Range("Y27").Value = '2018-11-15'
$ws.Range("Z27").NumberFormat = "@"
$ws.Range("Z27").Value = '00:00'
$ws.Range("AA27").NumberFormat = "@"
$ws.Range("AA27").Value = '2018-11-15'
$ws.Range("AB27").NumberFormat = "@"
$ws.Range("AB27").Value = '00:00'
$ws.Range("AW27").Value = 'Per Gustafsson'
$ws.Range("AX27").Value = 'Per Gustafsson'

# Row 28 (source old row 20)
$ws.Range("A28").Value = 74387740
$ws.Range("B28").Value = 77506
$ws.Range("P28").Value = 'Jonsmyrsudden, Vrm'
$ws.Range("Q28").Value = 376872.8754500371
$ws.Range("R28").Value = 6700377.027040293
$ws.Range("Y28").NumberFormat = "@"
$ws.Range("Y28").Value = '2018-11-15'
$ws.Range("Z28").NumberFormat = "@"
$ws.Range("Z28").Value = '00:00'
$ws.Range("AA28").NumberFormat = "@"
$ws.Range("AA28").Value = '2018-11-15'
$ws.Range("AB28").NumberFormat = "@"
$ws.Range("AB28").Value = '00:00'
$ws.Range("AW28").Value = 'Per Gustafsson'
$ws.Range("AX28").Value = 'Per Gustafsson'

# Row 29 (source old row 37)
$ws.Range("A29").Value = 74387734
$ws.Range("Q29").Value = 376901.8853666277
$ws.Range("R29").Value = 6699846.171224297

# Row 30 (source old row 38)
$ws.Range("A30").Value = 74387745
$ws.Range("Q30").Value = 376869.0183758087
$ws.Range("R30").Value = 6700030.180326126

# Row 31 (source old row 39)
$ws.Range("A31").Value = 74387746
$ws.Range("Q31").Value = 376867.8113538623
$ws.Range("R31").Value = 6700023.795488692

# Row 32 (source old row 21)
$ws.Range("A32").Value = 74387744
$ws.Range("B32").Value = 78098
$ws.Range("E32").Value = 6453
$ws.Range("F32").Value = 'Vedskivlav'
$ws.Range("G32").Value = 'Hertelidea botryosa'
$ws.Range("H32").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q32").Value = 376869.1209517682
$ws.Range("R32").Value = 6700412.245894314

# Row 33 (source old row 40)
$ws.Range("A33").Value = 74387741
$ws.Range("Q33").Value = 376871.1388198495
$ws.Range("R33").Value = 6699903.072625135

# Row 34 (source old row 41)
$ws.Range("A34").Value = 74387751
$ws.Range("Q34").Value = 376864.1398596731
$ws.Range("R34").Value = 6699973.996064176

# Row 35 (source old row 42)
$ws.Range("A35").Value = 74387732
$ws.Range("Q35").Value = 376913.9812352865
$ws.Range("R35").Value = 6699837.851543591

# Row 36 (source old row 43)
$ws.Range("A36").Value = 74387737
$ws.Range("Q36").Value = 376874.0559144435
$ws.Range("R36").Value = 6700018.146202448

# Row 37 (source old row 44)
$ws.Range("A37").Value = 74387753
$ws.Range("Q37").Value = 376858.7918946482
$ws.Range("R37").Value = 6699933.150422234

# Row 38 (source old row 50)
$ws.Range("A38").Value = 112322604
$ws.Range("B38").Value = 77650
$ws.Range("P38").Value = 'Gräsviggen, Vrm'
$ws.Range("Q38").Value = 377037
$ws.Range("R38").Value = 6699915
$ws.Range("Y38").NumberFormat = "@"
$ws.Range("Y38").Value = '2023-09-25'
$ws.Range("Z38").ClearContents() | Out-Null
$ws.Range("AA38").NumberFormat = "@"
$ws.Range("AA38").Value = '2023-09-25'
$ws.Range("AB38").ClearContents() | Out-Null
$ws.Range("AW38").Value = 'Helena Malmestrand'
$ws.Range("AX38").Value = 'Helena Malmestrand'

# Row 39 (source old row 24)
$ws.Range("A39").Value = 112323406
$ws.Range("B39").Value = 77650
$ws.Range("P39").Value = 'Jonsmyren, Vrm'
$ws.Range("Q39").Value = 376867
$ws.Range("R39").Value = 6700216
$ws.Range("Y39").NumberFormat = "@"
$ws.Range("Y39").Value = '2023-09-25'
$ws.Range("Z39").ClearContents() | Out-Null
$ws.Range("AA39").NumberFormat = "@"
$ws.Range("AA39").Value = '2023-09-25'
$ws.Range("AB39").ClearContents() | Out-Null
$ws.Range("AW39").Value = 'anders tedeholm'
$ws.Range("AX39").Value = 'anders tedeholm'

# Row 40 (source old row 27)
$ws.Range("A40").Value = 112323408
$ws.Range("B40").Value = 90837
$ws.Range("E40").Value = 5966
$ws.Range("F40").Value = 'Motaggsvamp'
$ws.Range("G40").Value = 'Sarcodon squamosus'
$ws.Range("H40").Value = '(Schaeff.) Quél.'
$ws.Range("P40").Value = 'Jonsmyren, Vrm'
$ws.Range("Q40").Value = 376863
$ws.Range("R40").Value = 6700246
$ws.Range("Y40").NumberFormat = "@"
$ws.Range("Y40").Value = '2023-09-25'
$ws.Range("Z40").ClearContents() | Out-Null
$ws.Range("AA40").NumberFormat = "@"
$ws.Range("AA40").Value = '2023-09-25'
$ws.Range("AB40").ClearContents() | Out-Null
$ws.Range("AW40").Value = 'anders tedeholm'
$ws.Range("AX40").Value = 'anders tedeholm'

# Row 41 (source old row 45)
$ws.Range("A41").Value = 112323397
$ws.Range("B41").Value = 77650
$ws.Range("P41").Value = 'Jonsmyren, Vrm'
$ws.Range("Q41").Value = 377010
$ws.Range("R41").Value = 6699884
$ws.Range("Y41").NumberFormat = "@"
$ws.Range("Y41").Value = '2023-09-25'
$ws.Range("Z41").ClearContents() | Out-Null
$ws.Range("AA41").NumberFormat = "@"
$ws.Range("AA41").Value = '2023-09-25'
$ws.Range("AB41").ClearContents() | Out-Null
$ws.Range("AW41").Value = 'anders tedeholm'
$ws.Range("AX41").Value = 'anders tedeholm'

# Row 42 (source old row 23)
$ws.Range("A42").Value = 112323404
$ws.Range("B42").Value = 78242
$ws.Range("E42").Value = 6453
$ws.Range("F42").Value = 'Vedskivlav'
$ws.Range("G42").Value = 'Hertelidea botryosa'
$ws.Range("H42").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("P42").Value = 'Jonsmyren, Vrm'
$ws.Range("Q42").Value = 376956
$ws.Range("R42").Value = 6700282
$ws.Range("Y42").NumberFormat = "@"
$ws.Range("Y42").Value = '2023-09-25'
$ws.Range("Z42").ClearContents() | Out-Null
$ws.Range("AA42").NumberFormat = "@"
$ws.Range("AA42").Value = '2023-09-25'
$ws.Range("AB42").ClearContents() | Out-Null
$ws.Range("AW42").Value = 'anders tedeholm'
$ws.Range("AX42").Value = 'anders tedeholm'

# Row 43 (source old row 46)
$ws.Range("A43").Value = 112323396
$ws.Range("B43").Value = 77650
$ws.Range("P43").Value = 'Jonsmyren, Vrm'
$ws.Range("Q43").Value = 376970
$ws.Range("R43").Value = 6699876
$ws.Range("Y43").NumberFormat = "@"
$ws.Range("Y43").Value = '2023-09-25'
$ws.Range("Z43").ClearContents() | Out-Null
$ws.Range("AA43").NumberFormat = "@"
$ws.Range("AA43").Value = '2023-09-25'
$ws.Range("AB43").ClearContents() | Out-Null
$ws.Range("AW43").Value = 'anders tedeholm'
$ws.Range("AX43").Value = 'anders tedeholm'

# Row 44 (source old row 53)
$ws.Range("A44").Value = 112323399
$ws.Range("B44").Value = 77650
$ws.Range("P44").Value = 'Jonsmyren, Vrm'
$ws.Range("Q44").Value = 377047
$ws.Range("R44").Value = 6699906
$ws.Range("Y44").NumberFormat = "@"
$ws.Range("Y44").Value = '2023-09-25'
$ws.Range("Z44").ClearContents() | Out-Null
$ws.Range("AA44").NumberFormat = "@"
$ws.Range("AA44").Value = '2023-09-25'
$ws.Range("AB44").ClearContents() | Out-Null
$ws.Range("AW44").Value = 'anders tedeholm'
$ws.Range("AX44").Value = 'anders tedeholm'

# Row 45 (source old row 55)
$ws.Range("A45").Value = 112322603
$ws.Range("B45").Value = 78242
$ws.Range("E45").Value = 6453
$ws.Range("F45").Value = 'Vedskivlav'
$ws.Range("G45").Value = 'Hertelidea botryosa'
$ws.Range("H45").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("P45").Value = 'Gräsviggen, Vrm'
$ws.Range("Q45").Value = 376933
$ws.Range("R45").Value = 6700074
$ws.Range("AW45").Value = 'Helena Malmestrand'
$ws.Range("AX45").Value = 'Helena Malmestrand'

# Row 46 (source old row 49)
$ws.Range("A46").Value = 112322635
$ws.Range("B46").Value = 78242
$ws.Range("E46").Value = 6453
$ws.Range("F46").Value = 'Vedskivlav'
$ws.Range("G46").Value = 'Hertelidea botryosa'
$ws.Range("H46").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("P46").Value = 'Gräsviggen, Vrm'
$ws.Range("Q46").Value = 376947
$ws.Range("R46").Value = 6699898
$ws.Range("AW46").Value = 'Helena Malmestrand'
$ws.Range("AX46").Value = 'Helena Malmestrand'

# Row 47 (source old row 28)
$ws.Range("A47").Value = 112323403
$ws.Range("P47").Value = 'Jonsmyren, Vrm'
$ws.Range("Q47").Value = 376976
$ws.Range("R47").Value = 6700256
$ws.Range("AW47").Value = 'anders tedeholm'
$ws.Range("AX47").Value = 'anders tedeholm'

# Row 48 (source old row 52)
$ws.Range("A48").Value = 112323400
$ws.Range("B48").Value = 77650
$ws.Range("D48").Value = 'NT'
$ws.Range("E48").Value = 6425
$ws.Range("F48").Value = 'Garnlav'
$ws.Range("G48").Value = 'Alectoria sarmentosa'
$ws.Range("H48").Value = '(Ach.) Ach.'
$ws.Range("Q48").Value = 376916
$ws.Range("R48").Value = 6699968

# Row 49 (source old row 47)
$ws.Range("A49").Value = 112322639
$ws.Range("B49").Value = 77650
$ws.Range("E49").Value = 6425
$ws.Range("F49").Value = 'Garnlav'
$ws.Range("G49").Value = 'Alectoria sarmentosa'
$ws.Range("H49").Value = '(Ach.) Ach.'
$ws.Range("Q49").Value = 376970
$ws.Range("R49").Value = 6699950

# Row 50 (source old row 26)
$ws.Range("A50").Value = 112322582
$ws.Range("Q50").Value = 376995
$ws.Range("R50").Value = 6700394

# Row 51 (source old row 56)
$ws.Range("A51").Value = 112322577
$ws.Range("B51").Value = 77403
$ws.Range("E51").Value = 228912
$ws.Range("F51").Value = 'Mörk kolflarnlav'
$ws.Range("G51").Value = 'Carbonicola myrmecina'
$ws.Range("H51").Value = '(Ach.) Bendiksby & Timdal'
$ws.Range("P51").Value = 'Gräsviggen, Vrm'
$ws.Range("Q51").Value = 376932
$ws.Range("R51").Value = 6700074
$ws.Range("AW51").Value = 'Helena Malmestrand'
$ws.Range("AX51").Value = 'Helena Malmestrand'

# Row 52 (source old row 54)
$ws.Range("A52").Value = 112323398
$ws.Range("Q52").Value = 377056
$ws.Range("R52").Value = 6699887

# Row 53 (source old row 25)
$ws.Range("A53").Value = 112323409
$ws.Range("B53").Value = 90837
$ws.Range("E53").Value = 5966
$ws.Range("F53").Value = 'Motaggsvamp'
$ws.Range("G53").Value = 'Sarcodon squamosus'
$ws.Range("H53").Value = '(Schaeff.) Quél.'
$ws.Range("Q53").Value = 376877
$ws.Range("R53").Value = 6700280

# Row 54 (source old row 51)
$ws.Range("A54").Value = 112323401
$ws.Range("B54").Value = 78242
$ws.Range("E54").Value = 6453
$ws.Range("F54").Value = 'Vedskivlav'
$ws.Range("G54").Value = 'Hertelidea botryosa'
$ws.Range("H54").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q54").Value = 376945
$ws.Range("R54").Value = 6700094

# Row 55 (source old row 48)
$ws.Range("A55").Value = 112323395
$ws.Range("B55").Value = 90814
$ws.Range("D55").Value = 'LC'
$ws.Range("E55").Value = 4364
$ws.Range("F55").Value = 'Dropptaggsvamp'
$ws.Range("G55").Value = 'Hydnellum ferrugineum'
$ws.Range("H55").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("P55").Value = 'Jonsmyren, Vrm'
$ws.Range("Q55").Value = 376769
$ws.Range("R55").Value = 6699865
$ws.Range("AW55").Value = 'anders tedeholm'
$ws.Range("AX55").Value = 'anders tedeholm'

# Row 56 (source old row 22)
$ws.Range("A56").Value = 112323405
$ws.Range("B56").Value = 77650
$ws.Range("E56").Value = 6425
$ws.Range("F56").Value = 'Garnlav'
$ws.Range("G56").Value = 'Alectoria sarmentosa'
$ws.Range("H56").Value = '(Ach.) Ach.'
$ws.Range("P56").Value = 'Jonsmyren, Vrm'
$ws.Range("Q56").Value = 376903
$ws.Range("R56").Value = 6700268
$ws.Range("AW56").Value = 'anders tedeholm'
$ws.Range("AX56").Value = 'anders tedeholm'
